# The sheet gained a new price-record row. A new row was inserted at
# position 39 (pushing the former rows 39-144 down to 40-145), and the
# new row 39 was populated with a fresh Papaya price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 39; existing rows 39..144 shift to 40..145
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new data record
$ws.Cells.Item(39, 1).Value = 10
$ws.Cells.Item(39, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(39, 3).Value = "La Araucanía"
$ws.Cells.Item(39, 4).Value = 45260
$ws.Cells.Item(39, 5).Value = 9
$ws.Cells.Item(39, 6).Value = "Fruta"
$ws.Cells.Item(39, 7).Value = 100108
$ws.Cells.Item(39, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(39, 9).Value = 100108004
$ws.Cells.Item(39, 10).Value = "Papaya"
$ws.Cells.Item(39, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(39, 12).Value = "Primera"
$ws.Cells.Item(39, 13).Value = 80
$ws.Cells.Item(39, 14).Value = 37500
$ws.Cells.Item(39, 15).Value = 37500
$ws.Cells.Item(39, 16).Value = 37500
$ws.Cells.Item(39, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(39, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(39, 19).Value = 2500
$ws.Cells.Item(39, 20).Value = 15
